$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Included Mixed Logit Model" — appends 7 more benders-cuts log rows
# (iterations alternating feasibility / optimality cuts) produced by a later
# solver run, including a widened dual-values dict and a widened
# problematic-ods list that now also covers the Baselwolf<->Chiasso pair.
# ---------------------------------------------------------------------------

# Text reused from rows already in the sheet (same values as before the edit)
$oldDuals = '{(''Chiasso'', ''Aarau''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Aarau'', ''Visp''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Visp'', ''Aarau''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}}'
$oldOd    = '[["Aarau", "Chiasso"], ["Chiasso", "Aarau"]]'

# New text values introduced by this edit
$newDuals = '{(''Aarau'', ''Chiasso''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Chiasso'', ''Aarau''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Aarau'', ''Visp''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Visp'', ''Aarau''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Baselwolf'', ''Chiasso''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}, (''Chiasso'', ''Baselwolf''): {''demand_dual'': 1.0, ''capacity_dual'': 1.0}}'
$newOd    = '[["Aarau", "Chiasso"], ["Chiasso", "Aarau"], ["Chiasso", "Baselwolf"]]'

$ts1 = '2025-11-28T03:02:04.062878'
$ts2 = '2025-11-28T03:02:06.371159'
$ts3 = '2025-11-28T03:05:45.389877'
$ts4 = '2025-11-28T03:05:48.897434'
$ts5 = '2025-11-28T03:09:09.108533'
$ts6 = '2025-11-28T03:09:12.880625'
$ts7 = '2025-11-28T03:09:17.193899'

# Seed the brand-new strings in the order they are first needed so the
# workbook's shared-string table grows the same way the original log export
# built it (new duals dict, then the seven timestamps, then the new od list).
$ws.Range("D14").Value = $newDuals
$ws.Range("E11").Value = $ts1
$ws.Range("E12").Value = $ts2
$ws.Range("E13").Value = $ts3
$ws.Range("E14").Value = $ts4
$ws.Range("E15").Value = $ts5
$ws.Range("E16").Value = $ts6
$ws.Range("E17").Value = $ts7
$ws.Range("F15").Value = $newOd

# Row 11 - feasibility cut
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "feasibility"
$ws.Range("E11").Value = $ts1
$ws.Range("F11").Value = $oldOd
$ws.Range("G11").Value = 233.1

# Row 12 - optimality cut
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "optimality"
$ws.Range("C12").Value = 38234.46124547702
$ws.Range("D12").Value = $oldDuals
$ws.Range("E12").Value = $ts2

# Row 13 - feasibility cut
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "feasibility"
$ws.Range("E13").Value = $ts3
$ws.Range("F13").Value = $oldOd
$ws.Range("G13").Value = 233.1

# Row 14 - optimality cut (new, wider duals dict + higher subproblem cost)
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "optimality"
$ws.Range("C14").Value = 65756.07328262705
$ws.Range("D14").Value = $newDuals
$ws.Range("E14").Value = $ts4

# Row 15 - feasibility cut (new, wider problematic-ods list + new rhs)
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "feasibility"
$ws.Range("E15").Value = $ts5
$ws.Range("F15").Value = $newOd
$ws.Range("G15").Value = 336.9240000000007

# Row 16 - feasibility cut (new, wider problematic-ods list + new rhs)
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "feasibility"
$ws.Range("E16").Value = $ts6
$ws.Range("F16").Value = $newOd
$ws.Range("G16").Value = 277.3889999999924

# Row 17 - optimality cut (new, wider duals dict + higher subproblem cost)
$ws.Range("A17").Value = 3
$ws.Range("B17").Value = "optimality"
$ws.Range("C17").Value = 71067.76083747702
$ws.Range("D17").Value = $newDuals
$ws.Range("E17").Value = $ts7
